$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.426.24"
$ws.Range("E2").Value = "  -2.89%  "
$ws.Range("D3").Value = "'1.986.28"
$ws.Range("E3").Value = "  -3.46%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'245.67"
$ws.Range("E5").Value = "  -3.00%  "
$ws.Range("E6").Value = "  -3.46%  "
$ws.Range("D7").Value = "'59.59"
$ws.Range("E7").Value = "  -12.51%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("D10").Value = "'57.23"
$ws.Range("E10").Value = "  -4.10%  "
$ws.Range("D11").Value = "'0.0825"
$ws.Range("E11").Value = "  +7.33%  "
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").Value = "'23.66"
$ws.Range("E13").Value = "  +4.50%  "
$ws.Range("D14").Value = "'0.867"
$ws.Range("E14").Value = "  -7.00%  "
$ws.Range("D15").Value = "'14.04"
$ws.Range("E15").Value = "  -5.85%  "
$ws.Range("D16").Value = "'2.276.69"
$ws.Range("E16").Value = "  -3.45%  "
$ws.Range("D17").Value = "'5.48"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "'1.980.54"
$ws.Range("E18").Value = "  -3.72%  "
$ws.Range("D19").Value = "'36.306.94"
$ws.Range("E19").Value = "  -2.88%  "
$ws.Range("D20").Value = "'70.51"
$ws.Range("E20").Value = "  -4.38%  "
$ws.Range("D21").Value = "'0.0₃0877"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'5.33"
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").Value = "'234.46"
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'2.58"
$ws.Range("E25").Value = "  -4.81%  "
$ws.Range("D26").Value = "'2.30"
$ws.Range("E26").Value = "  -5.28%  "
$ws.Range("D27").Value = "'10.00"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").Value = "'162.35"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "'19.91"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("E32").Value = "  -2.02%  "
$ws.Range("D33").Value = "'4.92"
$ws.Range("E33").Value = "  -6.13%  "
$ws.Range("E34").Value = "  +3.89%  "
$ws.Range("D35").Value = "'4.42"
$ws.Range("E35").Value = "  -5.52%  "
$ws.Range("D36").Value = "'6.24"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'2.27"
$ws.Range("E38").Value = "  -7.14%  "
$ws.Range("D39").Value = "'1.78"
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("D40").Value = "'2.97"
$ws.Range("E40").Value = "  -5.68%  "
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("D42").Value = "'0.0972"
$ws.Range("E42").Value = "  -6.30%  "
$ws.Range("E43").Value = "  -4.74%  "
$ws.Range("D44").Value = "'0.0214"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").Value = "'1.09"
$ws.Range("E45").Value = "  -4.96%  "
$ws.Range("D46").Value = "'16.29"
$ws.Range("E46").Value = "  -6.84%  "
$ws.Range("D47").Value = "'92.67"
$ws.Range("E47").Value = "  -4.88%  "
$ws.Range("D48").Value = "'7.50"
$ws.Range("E48").Value = "  -5.65%  "
$ws.Range("D49").Value = "'1.367.02"
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("E50").Value = "  -4.15%  "
$ws.Range("E51").Value = "  -2.35%  "
